$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '310.46'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '2.16%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '35.45'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-2.03%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.083'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.96%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08160'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.070'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-2.62%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.944'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.17%'
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.859'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '3.87%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9264'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.55%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1057'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '11.06%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1923'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '3.88%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09166'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '5.64%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03652'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '1.18%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09914'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.07%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001431'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.28%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005738'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '2.06%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.478'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.47%'
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.129'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.32%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-2.77%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.100'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.59%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2213'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-1.63%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04545'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.53%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001226'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.61%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004782'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.27%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-3.94%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0004450'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-6.36%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01970'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04890'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '3.92%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007570'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-2.79%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1380'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-0.37%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '16.06%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002150'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-2.52%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01163'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '3.98%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006620'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '3.98%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000750'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.15%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '183.60'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '252.23%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001501'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-21.07%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002101'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.15%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002001'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.15%'
